# Allow "year class" to be used in the collections parser.
# Rename the "Collection" column header to "Year Class" and update the
# sample collection values in the Template sheet to include the year
# (e.g. "WS" -> "2021 WS", "FP" -> "2021 FP").

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Template")

$ws.Range("G3").Value = "Year Class"
$ws.Range("G4").Value = "2021 FP"
$ws.Range("G5").Value = "2021 WS"
$ws.Range("G6").Value = "2021 WS"

$ws.Activate()
$ws.Range("G5").Select()
